# Auto-generated edit script
# Removes footnote markers like [1]..[6] from Vaccine names, joins multi-line
# brand/vaccine labels into single-line text, and fixes two duplicate brand-name
# entries (Recombivax HB, FluLaval Quadrivalent) that had stray extra rows in the
# shared string table because of inconsistent line breaks.

$wb = $excel.ActiveWorkbook


# --- Sheet 1: Pediatric Vaccine ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'DTaP '
$ws.Range("A3").Value = 'DTaP '
$ws.Range("A4").Value = 'DTaP '
$ws.Range("A5").Value = 'DTaP-IPV '
$ws.Range("A6").Value = 'DTaP-IPV '
$ws.Range("A7").Value = 'DTaP-Hep B-IPV '
$ws.Range("A8").Value = 'DTaP-IP-HI '
$ws.Range("A9").Value = 'e-IPV '
$ws.Range("A10").Value = 'Hepatitis A Pediatric '
$ws.Range("A11").Value = 'Hepatitis A Pediatric '
$ws.Range("A12").Value = 'Hepatitis A Pediatric '
$ws.Range("A13").Value = 'Hepatitis A Pediatric '
$ws.Range("A14").Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range("A15").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("A16").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("A17").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("B17").Value = 'Recombivax HB'
$ws.Range("A18").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("B18").Value = 'Recombivax HB'
$ws.Range("A19").Value = 'Hib '
$ws.Range("A20").Value = 'Hib '
$ws.Range("A21").Value = 'HIBMENCY '
$ws.Range("A22").Value = 'Hib '
$ws.Range("A23").Value = 'HPV - Human Papillomavirus 9-valent '
$ws.Range("A24").Value = 'MENB - Meningococcal Group B '
$ws.Range("A25").Value = 'MENB - Meningococcal Group B '
$ws.Range("A26").Value = 'MENB - Meningococcal Group B '
$ws.Range("A27").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A28").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A29").Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range("A30").Value = 'MMR/Varicella '
$ws.Range("A31").Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range("A33").Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range("A34").Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range("A35").Value = 'Rotavirus, Live, Oral, Oral '
$ws.Range("A36").Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A37").Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A38").Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A39").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A40").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A41").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A42").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A43").Value = 'Varicella '

# --- Sheet 2: Adult Vaccine ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Hepatitis A-Adult '
$ws.Range("A3").Value = 'Hepatitis A-Adult '
$ws.Range("A4").Value = 'Hepatitis A Adult '
$ws.Range("A5").Value = 'Hepatitis A Adult '
$ws.Range("A6").Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range("A7").Value = 'Hepatitis B-Adult '
$ws.Range("A8").Value = 'Hepatitis B-Adult '
$ws.Range("A9").Value = 'Hepatitis B-Adult '
$ws.Range("A10").Value = 'Hepatitis B-Adult '
$ws.Range("A11").Value = 'HPV-Human Papillomavirus 9 Valent '
$ws.Range("A12").Value = 'Measles, Mumps,  Rubella-Adult '
$ws.Range("A13").Value = 'Meningococcal Conjugate '
$ws.Range("A14").Value = 'Meningococcal Conjugate '
$ws.Range("A15").Value = 'MENB - Meningococcal Group B '
$ws.Range("A16").Value = 'MENB - Meningococcal Group B '
$ws.Range("A17").Value = 'MENB - Meningococcal Group B '
$ws.Range("A18").Value = 'Pneumococcal 13-valent  (Adult)'
$ws.Range("A21").Value = 'Tetanus and Diphtheria Toxoids'
$ws.Range("G21").Value = ''
$ws.Range("I21").Value = ''
$ws.Range("A22").Value = 'Tetanus and Diphtheria Toxoids'
$ws.Range("G22").Value = ''
$ws.Range("I22").Value = ''
$ws.Range("A23").Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A24").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A25").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A26").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A27").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A28").Value = 'Varicella-Adult '

# --- Sheet 3: Pediatric Influenza Vaccine ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Influenza  (Age 6 months and older)'
$ws.Range("B2").Value = 'Fluzone Quadrivalent'
$ws.Range("A3").Value = 'Influenza  (Age 6-35 months)'
$ws.Range("B3").Value = 'Fluzone Quadrivalent Pediatric dose No Preservative'
$ws.Range("A4").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B4").Value = 'Fluzone Quadrivalent No-Preservative'
$ws.Range("A5").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B5").Value = 'Fluzone Quadrivalent No-Preservative'
$ws.Range("A6").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B6").Value = 'Fluarix Quadrivalent Preservative Free'
$ws.Range("A7").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B7").Value = 'FluLaval Quadrivalent'
$ws.Range("A8").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B8").Value = 'FluLaval Quadrivalent'
$ws.Range("A9").Value = 'Influenza  (Age 4 years and older)'
$ws.Range("A10").Value = 'Influenza  (Age 9 years and older)'
$ws.Range("B10").Value = 'Afluria No Preservative'
$ws.Range("D10").Value = '10 pack-1 dose syringe'
$ws.Range("A11").Value = 'Influenza  (Age 9 years and older)'
$ws.Range("A12").Value = 'Influenza  (Age 4 years and older)'

# --- Sheet 4: Adult Influenza Vaccine ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 'Influenza  (Age 6 months and older)'
$ws.Range("B2").Value = 'Fluzone Quadrivalent'
$ws.Range("A3").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B3").Value = 'Fluzone Quadrivalent No Preservative'
$ws.Range("A4").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B4").Value = 'Fluzone Quadrivalent No Preservative'
$ws.Range("A5").Value = 'Influenza  (Age 4 years and older)'
$ws.Range("A6").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("B6").Value = 'Fluarix Quadrivalent Preservative Free'
$ws.Range("A7").Value = 'Influenza  (Age 36 months and older)'
$ws.Range("A8").Value = 'Influenza  (Age 9 years and older)'
$ws.Range("B8").Value = 'Afluria No Preservative'
$ws.Range("D8").Value = '10 pack-1 dose syringe'
$ws.Range("A9").Value = 'Influenza  (Age 9 years and older)'
$ws.Range("A10").Value = 'Influenza  (Age 18 years and older)'
